$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "What opportunities exist for Monalco Mining..." problem-statement
# textbox by its distinctive text rather than a hard-coded shape index.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -like "*Mining to reduce ore-crusher maintenance costs by 20%*") {
            $target = $sh
        }
    }
}

$tr = $target.TextFrame.TextRange
$full = $tr.Text

# The sentence originally ended "...(from $45M to $36M in 2020?" - missing the
# closing parenthesis after $36M. Fix it to "...$36M) in 2020?" while keeping
# the untouched "$36M " characters re-typed as their own run (matching how
# PowerPoint splits a run when the text inside it is edited in place).
$pos = $full.IndexOf("`$36M in 2020?")
if ($pos -ge 0) {
    $sub = $tr.Characters($pos + 1, 5)
    $sub.Text = "`$36M) "
}
